# Auto-generated: apply crypto price/volume update diff to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.820.27"

# Row 3
$ws.Range("D3").Value = "2.317.30"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.61%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.73"
$ws.Range("E5").Value = "  -1.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.70"
$ws.Range("E6").Value = "  +4.68%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.26"
$ws.Range("E10").Value = "  +2.14%  "

# Row 11
$ws.Range("E11").Value = "  +0.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.40"
$ws.Range("E12").Value = "  -0.95%  "

# Row 13
$ws.Range("E13").Value = "  -0.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.997"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16
$ws.Range("D16").Value = "2.665.75"
$ws.Range("E16").Value = "  +0.70%  "

# Row 17
$ws.Range("D17").Value = "2.313.59"
$ws.Range("E17").Value = "  +0.78%  "

# Row 18
$ws.Range("D18").Value = "43.004.47"
$ws.Range("E18").Value = "  +1.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.51"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("E21").Value = "  -4.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.78"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -1.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.67"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.39%  "

# Row 27
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.79"
$ws.Range("E27").Value = "  +17.92%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.00"
$ws.Range("E28").Value = "  +0.84%  "

# Row 29
$ws.Range("E29").Value = "  +1.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.70"
$ws.Range("E30").Value = "  +4.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.42"
$ws.Range("E31").Value = "  -0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.14"
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0870"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").Value = "  +6.72%  "

# Row 35
$ws.Range("E35").Value = "  -0.65%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("E37").Value = "  +2.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("E38").Value = "  +1.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("E39").Value = "  +4.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("E40").Value = "  -1.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.60"
$ws.Range("E41").Value = "  +0.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.87"
$ws.Range("E42").Value = "  +10.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.17"
$ws.Range("E43").Value = "  +1.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.233"
$ws.Range("E44").Value = "  +2.13%  "

# Row 45
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.41"
$ws.Range("E46").Value = "  +0.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.90"
$ws.Range("E47").Value = "  -1.88%  "

# Row 48
$ws.Range("D48").Value = "1.685.85"
$ws.Range("E48").Value = "  +1.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.19"
$ws.Range("E49").Value = "  -3.79%  "

# Row 50
$ws.Range("E50").Value = "  +0.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.25"
$ws.Range("E51").Value = "  -0.12%  "

